$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

# Column A holds a plain date-like text string in this sheet (not a real date
# value), so force text formatting to stop Excel's automatic date detection,
# then clear the format again so no stray style survives on the cell.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/18/2026"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = 12566.24
$ws.Cells.Item($row, 3).Value = 0.2341426964151604
$ws.Cells.Item($row, 4).Value = 0.7658573035848396
$ws.Cells.Item($row, 5).Value = -152.55
$ws.Cells.Item($row, 6).Value = -22.27
$ws.Cells.Item($row, 7).Value = -21117.81
$ws.Cells.Item($row, 8).Value = -68.69
$ws.Cells.Item($row, 9).Value = -260.56
$ws.Cells.Item($row, 10).Value = -8.140000000000001
$ws.Cells.Item($row, 11).Value = -21378.37
$ws.Cells.Item($row, 12).Value = -62.98
